{"js": "// The diff removes the custom paragraph style \"Abstract Title\" (styleId\n// \"AbstractTitle\") from styles.xml, and changes the \"Abstract\" style's\n// paragraph spacing \"before\" value from 100 twips (5pt) to 300 twips (15pt).\n// Neither style is referenced by any paragraph in the document body, so this\n// is purely a style-definition-level edit.\n\nconst styles = context.document.getStyles();\n\n// Remove the \"Abstract Title\" style entirely.\nconst abstractTitleStyle = styles.getByNameOrNullObject(\"Abstract Title\");\nabstractTitleStyle.load(\"isNullObject\");\nawait context.sync();\n\nif (!abstractTitleStyle.isNullObject) {\n    abstractTitleStyle.delete();\n    await context.sync();\n}\n\n// Update the \"Abstract\" style's space-before from 5pt (100 twips) to 15pt\n// (300 twips) to match its space-after.\nconst abstractStyle = styles.getByNameOrNullObject(\"Abstract\");\nabstractStyle.load(\"isNullObject\");\nawait context.sync();\n\nif (!abstractStyle.isNullObject) {\n    abstractStyle.paragraphFormat.spaceBefore = 15; // points; 15pt == 300 twips\n    await context.sync();\n}\n", "ps1": "# The diff removes the custom paragraph style \"Abstract Title\" (styleId\n# \"AbstractTitle\") from styles.xml, and changes the \"Abstract\" style's\n# paragraph spacing \"before\" value from 100 twips (5pt) to 300 twips (15pt).\n# Neither style is referenced by any paragraph in the document body, so this\n# is purely a style-definition-level edit.\n\n$d = $word.ActiveDocument\n\n# Remove the \"Abstract Title\" style entirely.\ntry {\n    $abstractTitleStyle = $d.Styles.Item(\"Abstract Title\")\n    if ($abstractTitleStyle -ne $null) {\n        $abstractTitleStyle.Delete()\n    }\n} catch {\n    # Style already absent - nothing to remove.\n}\n\n# Update the \"Abstract\" style's space-before from 5pt (100 twips) to 15pt\n# (300 twips) to match its space-after.\ntry {\n    $abstractStyle = $d.Styles.Item(\"Abstract\")\n    if ($abstractStyle -ne $null) {\n        $abstractStyle.ParagraphFormat.SpaceBefore = 15\n    }\n} catch {\n    # Style missing - nothing to update.\n}\n"}
